# Highlight specific bullet points in green (00FF00) on slides 7 and 8,
# matching the "void Matrix2DIdentity" Step 0 / Step 1 checklist items.

$p = $ppt.ActivePresentation

# Green highlight color, expressed the VBA RGB() way (0x00BBGGRR -> here 0x00FF00 = RGB(0,255,0))
$green = 65280

function Set-ParaHighlight($textRange, [int]$paraIndex) {
    $para = $textRange.Paragraphs($paraIndex, 1)
    $para.Font.Highlight.RGB = $green
}

# Slide 7 ("Step 0") - Content Placeholder 2
$slide7 = $p.Slides.Item(7)
$tr7 = $slide7.Shapes.Item(2).TextFrame.TextRange
Set-ParaHighlight $tr7 1   # "Unzip the Project 3 materials into a clean folder"
Set-ParaHighlight $tr7 2   # "Add and commit the files to version control"

# Slide 8 ("Step 1") - Content Placeholder 2
$slide8 = $p.Slides.Item(8)
$tr8 = $slide8.Shapes.Item(2).TextFrame.TextRange
Set-ParaHighlight $tr8 1   # "Integrate code from Project 2 into Project 3"
Set-ParaHighlight $tr8 2   # "Create new .c modules and add stub functions"
Set-ParaHighlight $tr8 4   # "Add temporary return values and unreferenced parameter macros..."
Set-ParaHighlight $tr8 6   # "Add and commit the files to version control"
